$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Update client names (column B)
$ws1.Range("B2").Value = "Winn"
$ws1.Range("B3").Value = "Keevil"
$ws1.Range("B4").Value = "Howard"
$ws1.Range("B5").Value = "Markfield"
$ws1.Range("B6").Value = "Layne"

# Update Rate (column E) and Total (column F) for rows 2-6
$ws1.Range("E2:E6").Value = 90
$ws1.Range("F2:F6").Value = 720

# Update subtotal/total rows in column F
$ws1.Range("F8").Value = 3600
$ws1.Range("F11").Value = 3600
$ws1.Range("F13").Value = 3600

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Update client names (column D) to match sheet 1
$ws2.Range("D2").Value = "Winn"
$ws2.Range("D3").Value = "Keevil"
$ws2.Range("D4").Value = "Howard"
$ws2.Range("D5").Value = "Markfield"
$ws2.Range("D6").Value = "Layne"

# Update Rate (column F) and Total (column G) for rows 2-6
$ws2.Range("F2:F6").Value = 90
$ws2.Range("G2:G6").Value = 720

# Update Employee ID (column B) for rows 2-6
$ws2.Range("B2:B6").Value = "emp_75yd72zj"
